$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells that could be misread as numbers by Excel's auto-detection
# are written with a temporary Text number format, then the format is reset
# back to Normal style so no stray style survives in the saved file.
function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "26.895.60"
Set-TextValue "E2" "  -0.91%  "
Set-TextValue "D3" "1.868.93"
Set-TextValue "E3" "  +0.08%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "305.97"
Set-TextValue "E5" "  -0.04%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  +0.02%  "
Set-TextValue "D7" "0.5096"
Set-TextValue "E7" "  -1.02%  "
Set-TextValue "E8" "  -2.63%  "
Set-TextValue "D9" "0.07163"
Set-TextValue "E9" "  +0.26%  "
Set-TextValue "D10" "0.8898"
Set-TextValue "E10" "  -0.12%  "
Set-TextValue "D11" "20.56"
Set-TextValue "E11" "  -0.84%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.07484"
Set-TextValue "E12" "  -0.87%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.868.93"
Set-TextValue "E13" "  +0.19%  "
Set-TextValue "D14" "94.55"
Set-TextValue "E14" "  +5.53%  "
Set-TextValue "D15" "5.219"
Set-TextValue "E15" "  -1.71%  "
Set-TextValue "E16" "  +0.11%  "
Set-TextValue "D17" "0.000008484"
Set-TextValue "E17" "  -0.04%  "
Set-TextValue "E18" "  +0.16%  "
Set-TextValue "D19" "1.000"
Set-TextValue "E19" "  -0.01%  "
Set-TextValue "D20" "26.931.43"
Set-TextValue "E20" "  -0.87%  "
Set-TextValue "D21" "4.999"
Set-TextValue "E21" "  +0.09%  "
Set-TextValue "D22" "2.116.06"
Set-TextValue "D23" "10.33"
Set-TextValue "E23" "  -1.40%  "
Set-TextValue "D24" "6.368"
Set-TextValue "E24" "  -1.31%  "
Set-TextValue "D25" "147.70"
Set-TextValue "E25" "  +1.16%  "
Set-TextValue "D26" "1.779"
Set-TextValue "E26" "  -3.04%  "
Set-TextValue "E27" "  -0.55%  "
Set-TextValue "D28" "2.082"
Set-TextValue "E28" "  -0.33%  "
Set-TextValue "D29" "113.38"
Set-TextValue "E29" "  +0.57%  "
Set-TextValue "D30" "4.669"
Set-TextValue "E30" "  +0.19%  "
Set-TextValue "E31" "  +0.47%  "
Set-TextValue "D32" "0.09125"
Set-TextValue "E32" "  -1.09%  "
Set-TextValue "D33" "0.05027"
Set-TextValue "E33" "  -1.66%  "
Set-TextValue "D34" "2.990"
Set-TextValue "E34" "  -2.80%  "
Set-TextValue "D35" "0.7468"
Set-TextValue "E35" "  +3.03%  "
Set-TextValue "D36" "1.151"
Set-TextValue "E36" "  -0.78%  "
Set-TextValue "D37" "3.213"
Set-TextValue "E37" "  +3.74%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "2.508"
Set-TextValue "E38" "  +0.18%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D39" "0.5595"
Set-TextValue "E39" "  +5.74%  "
Set-TextValue "E40" "  -2.46%  "
Set-TextValue "D41" "1.072"
Set-TextValue "E41" "  -0.50%  "
Set-TextValue "D42" "6.584"
Set-TextValue "E42" "  +1.22%  "
Set-TextValue "D43" "115.41"
Set-TextValue "E43" "  -1.07%  "
Set-TextValue "D44" "8.540"
Set-TextValue "E44" "  +2.51%  "
Set-TextValue "D45" "0.1482"
Set-TextValue "E45" "  +0.84%  "
Set-TextValue "D46" "0.4767"
Set-TextValue "E46" "  +3.17%  "
Set-TextValue "D47" "1.000"
Set-TextValue "E47" "  +0.03%  "
Set-TextValue "D48" "10.11"
Set-TextValue "E48" "  +1.44%  "
Set-TextValue "D49" "1.554"
Set-TextValue "E49" "  -0.31%  "
Set-TextValue "D50" "37.00"
Set-TextValue "E50" "  +0.95%  "
Set-TextValue "D51" "62.94"
Set-TextValue "E51" "  -1.04%  "
